$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-22 04:18:18"
$ws.Range("E3").Value = "2026-02-22 04:18:20"
$ws.Range("L3").Value = "17.3 km/h - 111º 3:38 TU"
$ws.Range("M3").Value = "4.1 °C 3:30 TU"
$ws.Range("N3").Value = "1.4 °C 3:52 TU"
$ws.Range("E4").Value = "2026-02-22 04:18:23"
$ws.Range("H4").Value = "'79%"
$ws.Range("J4").Value = "1028.3 hPa"
$ws.Range("K4").Value = "-0.1 MJ/m2"
$ws.Range("O4").Value = "6.8 °C"
$ws.Range("E5").Value = "2026-02-22 04:18:25"
$ws.Range("O5").Value = "4.9 °C"
$ws.Range("E6").Value = "2026-02-22 04:18:27"
$ws.Range("J6").Value = "1028.0 hPa"
$ws.Range("N6").Value = "6.3 °C 3:57 TU"
$ws.Range("O6").Value = "8.1 °C"
$ws.Range("E7").Value = "2026-02-22 04:18:30"
$ws.Range("J7").Value = "1028.2 hPa"
$ws.Range("O7").Value = "12.0 °C"
$ws.Range("E8").Value = "2026-02-22 04:18:32"
$ws.Range("J8").Value = "1028.1 hPa"
$ws.Range("O8").Value = "11.8 °C"
$ws.Range("E9").Value = "2026-02-22 04:18:34"
$ws.Range("H9").Value = "'95%"
$ws.Range("N9").Value = "2.6 °C 3:50 TU"
$ws.Range("O9").Value = "4.9 °C"
$ws.Range("E10").Value = "2026-02-22 04:18:37"
$ws.Range("L10").Value = "5.4 km/h - 75º 3:30 TU"
$ws.Range("E11").Value = "2026-02-22 04:18:39"
$ws.Range("N11").Value = "1.0 °C 3:40 TU"
$ws.Range("E12").Value = "2026-02-22 04:18:42"
$ws.Range("N12").Value = "3.1 °C 3:59 TU"
$ws.Range("O12").Value = "5.0 °C"
$ws.Range("E13").Value = "2026-02-22 04:18:44"
$ws.Range("H13").Value = "'91%"
$ws.Range("O13").Value = "-2.3 °C"
$ws.Range("E14").Value = "2026-02-22 04:18:46"
$ws.Range("E15").Value = "2026-02-22 04:18:49"
$ws.Range("H15").Value = "'85%"
$ws.Range("N15").Value = "2.4 °C 3:30 TU"
$ws.Range("O15").Value = "4.9 °C"
$ws.Range("E16").Value = "2026-02-22 04:18:51"
$ws.Range("H16").Value = "'24%"
$ws.Range("E17").Value = "2026-02-22 04:18:53"
$ws.Range("H17").Value = "'23%"
$ws.Range("E18").Value = "2026-02-22 04:18:56"
$ws.Range("J18").Value = "1028.7 hPa"
$ws.Range("N18").Value = "0.8 °C 3:55 TU"
$ws.Range("O18").Value = "1.6 °C"
$ws.Range("E19").Value = "2026-02-22 04:18:58"
$ws.Range("H19").Value = "'58%"
$ws.Range("O19").Value = "8.0 °C"
$ws.Range("E20").Value = "2026-02-22 04:19:01"
$ws.Range("H20").Value = "'41%"
$ws.Range("O20").Value = "0.4 °C"
$ws.Range("E21").Value = "2026-02-22 04:19:03"
$ws.Range("H21").Value = "'76%"
$ws.Range("J21").Value = "1033.4 hPa"
$ws.Range("N21").Value = "1.6 °C 3:32 TU"
$ws.Range("O21").Value = "3.0 °C"
$ws.Range("E22").Value = "2026-02-22 04:19:05"
$ws.Range("H22").Value = "'26%"
$ws.Range("L22").Value = "13.0 km/h - 325º 3:58 TU"
$ws.Range("O22").Value = "3.0 °C"
$ws.Range("E23").Value = "2026-02-22 04:19:08"
$ws.Range("L23").Value = "8.6 km/h - 2º 3:56 TU"
$ws.Range("E24").Value = "2026-02-22 04:19:10"
$ws.Range("J24").Value = "1031.6 hPa"
$ws.Range("L24").Value = "12.2 km/h - 161º 3:46 TU"
$ws.Range("M24").Value = "2.9 °C 3:48 TU"
$ws.Range("O24").Value = "1.1 °C"
$ws.Range("E25").Value = "2026-02-22 04:19:13"
$ws.Range("O25").Value = "4.4 °C"
$ws.Range("E26").Value = "2026-02-22 04:19:15"
$ws.Range("J26").Value = "1028.7 hPa"
$ws.Range("L26").Value = "11.9 km/h - 4º 3:35 TU"
$ws.Range("E27").Value = "2026-02-22 04:19:18"
$ws.Range("H27").Value = "'32%"
$ws.Range("N27").Value = "3.2 °C 3:45 TU"
$ws.Range("O27").Value = "3.8 °C"
$ws.Range("E28").Value = "2026-02-22 04:19:20"
$ws.Range("J28").Value = "1030.2 hPa"
$ws.Range("N28").Value = "1.0 °C 3:53 TU"
$ws.Range("O28").Value = "2.4 °C"
$ws.Range("E29").Value = "2026-02-22 04:19:22"
$ws.Range("H29").Value = "'94%"
$ws.Range("N29").Value = "2.7 °C 3:47 TU"
$ws.Range("O29").Value = "5.0 °C"
$ws.Range("E30").Value = "2026-02-22 04:19:25"
$ws.Range("J30").Value = "1028.0 hPa"
$ws.Range("O30").Value = "8.2 °C"
$ws.Range("E31").Value = "2026-02-22 04:19:27"
$ws.Range("H31").Value = "'60%"
$ws.Range("N31").Value = "10.4 °C 3:43 TU"
$ws.Range("O31").Value = "12.5 °C"
$ws.Range("E32").Value = "2026-02-22 04:19:29"
$ws.Range("E33").Value = "2026-02-22 04:19:32"
$ws.Range("H33").Value = "'64%"
$ws.Range("N33").Value = "0.1 °C 3:49 TU"
$ws.Range("O33").Value = "1.5 °C"
$ws.Range("E34").Value = "2026-02-22 04:19:34"
$ws.Range("H34").Value = "'45%"
$ws.Range("L34").Value = "18.0 km/h - 20º 3:47 TU"
$ws.Range("M34").Value = "5.0 °C 3:52 TU"
$ws.Range("O34").Value = "2.2 °C"
$ws.Range("E35").Value = "2026-02-22 04:19:37"
$ws.Range("H35").Value = "'37%"
$ws.Range("J35").Value = "1031.8 hPa"
$ws.Range("M35").Value = "6.8 °C 3:38 TU"
$ws.Range("O35").Value = "5.8 °C"
$ws.Range("E36").Value = "2026-02-22 04:19:39"
$ws.Range("J36").Value = "1027.9 hPa"
$ws.Range("O36").Value = "6.3 °C"
$ws.Range("E37").Value = "2026-02-22 04:19:41"
$ws.Range("J37").Value = "1034.0 hPa"
$ws.Range("N37").Value = "-0.9 °C 3:53 TU"
$ws.Range("O37").Value = "-0.2 °C"
$ws.Range("E38").Value = "2026-02-22 04:19:44"
$ws.Range("H38").Value = "'79%"
$ws.Range("O38").Value = "5.5 °C"
$ws.Range("E39").Value = "2026-02-22 04:19:46"
$ws.Range("L39").Value = "16.9 km/h - 247º 3:58 TU"
$ws.Range("N39").Value = "3.1 °C 3:56 TU"
$ws.Range("E40").Value = "2026-02-22 04:19:49"
$ws.Range("G40").Value = "2 cm"
$ws.Range("H40").Value = "'74%"
$ws.Range("J40").Value = "1033.0 hPa"
$ws.Range("E41").Value = "2026-02-22 04:19:51"
$ws.Range("J41").Value = "1028.7 hPa"
$ws.Range("O41").Value = "5.6 °C"
$ws.Range("E42").Value = "2026-02-22 04:19:53"
$ws.Range("H42").Value = "'98%"
$ws.Range("N42").Value = "3.6 °C 3:58 TU"
$ws.Range("O42").Value = "5.2 °C"
$ws.Range("E43").Value = "2026-02-22 04:19:55"
$ws.Range("N43").Value = "1.1 °C 3:59 TU"
$ws.Range("O43").Value = "2.4 °C"
$ws.Range("E44").Value = "2026-02-22 04:19:58"
$ws.Range("H44").Value = "'48%"
$ws.Range("O44").Value = "0.2 °C"
$ws.Range("E45").Value = "2026-02-22 04:20:00"
$ws.Range("J45").Value = "1031.7 hPa"
$ws.Range("E46").Value = "2026-02-22 04:20:02"
$ws.Range("J46").Value = "1031.4 hPa"
$ws.Range("N46").Value = "0.6 °C 3:49 TU"
$ws.Range("O46").Value = "1.6 °C"
